$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - subject max ROM values (columns B-E)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 - updated passive torque values (columns B-E)
$ws.Range("B2").Value = 49.443367045917803
$ws.Range("C2").Value = 21.928985080686441
$ws.Range("D2").Value = 56.186296641910367
$ws.Range("E2").Value = 20.821599490754718

# Row 3 - updated passive torque values (columns B-E)
$ws.Range("B3").Value = 46.028977094461943
$ws.Range("C3").Value = 21.561630270302388
$ws.Range("D3").Value = 50.973192801965183
$ws.Range("E3").Value = 17.426995298519632

# Update the selected range to reflect the new region of interest
$ws.Range("B1:E3").Select()
